# Applies the "Mappa di traduzione aggiornata" update to Foglio1.
# D6, D9, D19, D21, D22, D23 get new/changed text (the translation map was
# updated); column C is widened to fit the new, longer snippet in D6; and
# the sheet's selection is moved from D26 to C23 (closer to the edited area).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- D6: "click su Spinner" Robotium snippet -> longer multi-step snippet ---
$ws.Range("D6").Value = "Spinner s = (Spinner)solo.getView(R.id.idspinner);`n        solo.clickOnView(s);`n        solo.scrollToTop(); // I put this in here so that it always keeps the list at start`n        // select the 3th item in the spinner`n        solo.clickOnView(solo.getView(TextView.class, 3));"
$ws.Range("D6").WrapText = $true

# --- D9: previously empty -> same Robotium snippet used in D8 ---
$ws.Range("D9").Value = $ws.Range("D8").Value2

# --- D19: listview/"universita" snippet -> "//" ---
$ws.Range("D19").Value = "//"

# --- D21 / D22: previously empty -> "//" (same style family as D19/D17) ---
$ws.Range("D21").Value = "//"
$ws.Range("D21").Font.Name = "Arial"
$ws.Range("D21").Font.Size = 10
$ws.Range("D21").Font.Color = 2236962

$ws.Range("D22").Value = "//"
$ws.Range("D22").Font.Name = "Arial"
$ws.Range("D22").Font.Size = 10
$ws.Range("D22").Font.Color = 2236962

# --- D23: TextView "ciao!" snippet -> solo.enterText(...) snippet ---
$ws.Range("D23").Value = 'solo.enterText((EditText) solo.getView(R.id.idtextView), "idinputText!");'

# --- Column C: widen to fit the new content ---
$ws.Columns.Item(3).ColumnWidth = 73.3

# --- Move the visible selection to C23 (near the edited rows) ---
$ws.Range("C23").Select()
